# Scenario 8.xlsx edits
# - Decrease probability of fetal death before 4 weeks from conception
#   (Phase1!B2:B5 and Phase2!C2:C5, with their dependent "1 - ..." formulas
#   recalculating automatically).
# - Modify the Phase4 analysis formula for C5 from "=C3" to "=1*C3" so SEs
#   calculate correctly for truth.
# - Update the threaded/legacy comment on Phase4!C4 to reflect the new
#   treatment effect assumption.
# - Leave behind the new cell selections that resulted from the edits.

$wb = $excel.ActiveWorkbook

# --- Phase1: halve the probabilities in B2:B5 (D column is "=1-C-B", recalcs automatically)
$ws1 = $wb.Worksheets.Item("Phase1")
$ws1.Range("B2").Value = 0.1
$ws1.Range("B3").Value = 0.1
$ws1.Range("B4").Value = 0.05
$ws1.Range("B5").Value = 0.05

# --- Phase2: halve the probabilities in C2:C5 (E column is "=1-D-C", recalcs automatically)
$ws2 = $wb.Worksheets.Item("Phase2")
$ws2.Range("C2").Value = 0.1
$ws2.Range("C3").Value = 0.1
$ws2.Range("C4").Value = 0.05
$ws2.Range("C5").Value = 0.05

# --- Phase4: fix the C5 formula and update the comment text on C4
$ws4 = $wb.Worksheets.Item("Phase4")
$ws4.Range("C5").Formula = "=1*C3"

$commentThreaded = $ws4.Range("C4").CommentThreaded
[void]$commentThreaded.Text("Treatment DECREASES risk of SGA, RR = 0.8")

# --- Leave the selections where the author left them after editing each sheet
[void]$ws1.Range("B2:B5").Select()
[void]$ws2.Range("C2:C5").Select()

# Restore the originally active sheet/tab (Phase4) as the final state
[void]$ws4.Activate()
